# LMS8 upload and results
# Update the "< 3 Months" / "3-5 Months" survey-year tags in rows 2 and 12
# (F/H swap 2017->2022, I/K swap 2022->2017) on both sheets, then leave the
# second sheet ("wetight percentage by province") as the active tab with
# C2:L21 selected, matching the saved workbook state.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Range("F2").Value = 2022
    $ws.Range("H2").Value = 2022
    $ws.Range("I2").Value = 2017
    $ws.Range("K2").Value = 2017

    $ws.Range("F12").Value = 2022
    $ws.Range("H12").Value = 2022
    $ws.Range("I12").Value = 2017
    $ws.Range("K12").Value = 2017

    $ws.Range("C2:L21").Select() | Out-Null
}

$ws2 = $wb.Worksheets.Item(2)
$ws2.Activate() | Out-Null
$ws2.Range("C2:L21").Select() | Out-Null
